# Apply the updated cryptocurrency snapshot values (price + 1h volume change,
# plus the FraxShare/MXToken row swap) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores figures as text (e.g. "26.776.28", "1.01", "0.0619")
# in the source data. Plain decimal-looking strings would otherwise be
# auto-converted to numbers by Excel, so force text with a leading apostrophe
# (the same trick a user would use typing into the grid) whenever the new
# value parses as a plain number.
function Set-CellText {
    param($Cell, [string]$Text)
    if ($Text -match '^[-+]?[0-9]*\.?[0-9]+$') {
        $Cell.Value = "'" + $Text
    } else {
        $Cell.Value = $Text
    }
}

# Row 2: Bitcoin
Set-CellText $ws.Range("D2") '26.776.28'
Set-CellText $ws.Range("E2") '  +0.20%  '

# Row 3: Ethereum
Set-CellText $ws.Range("D3") '1.600.79'
Set-CellText $ws.Range("E3") '  +0.08%  '

# Row 4: TetherUSD
Set-CellText $ws.Range("D4") '1.01'
Set-CellText $ws.Range("E4") '  +0.20%  '

# Row 5: BNB
Set-CellText $ws.Range("D5") '211.90'
Set-CellText $ws.Range("E5") '  +0.26%  '

# Row 6: XRP
Set-CellText $ws.Range("E6") '  -0.18%  '

# Row 7: USDC
Set-CellText $ws.Range("D7") '1.01'
Set-CellText $ws.Range("E7") '  +0.20%  '

# Row 8: Dogecoin
Set-CellText $ws.Range("D8") '0.0619'
Set-CellText $ws.Range("E8") '  -0.02%  '

# Row 9: Cardano
Set-CellText $ws.Range("E9") '  -0.16%  '

# Row 10: Solana
Set-CellText $ws.Range("D10") '19.67'
Set-CellText $ws.Range("E10") '  +0.68%  '

# Row 11: TRON
Set-CellText $ws.Range("D11") '0.0847'
Set-CellText $ws.Range("E11") '  +0.78%  '

# Row 12: WrappedliquidstakedEther2.0
Set-CellText $ws.Range("D12") '1.829.77'
Set-CellText $ws.Range("E12") '  +0.33%  '

# Row 13: WrappedEther
Set-CellText $ws.Range("D13") '1.613.89'
Set-CellText $ws.Range("E13") '  +0.48%  '

# Row 14: Polkadot
Set-CellText $ws.Range("E14") '  +0.57%  '

# Row 15: Polygon
Set-CellText $ws.Range("D15") '0.526'
Set-CellText $ws.Range("E15") '  +0.51%  '

# Row 16: Litecoin
Set-CellText $ws.Range("D16") '65.19'
Set-CellText $ws.Range("E16") '  -0.40%  '

# Row 17: ShibaInu
Set-CellText $ws.Range("D17") '0.0₃0743'
Set-CellText $ws.Range("E17") '  -2.69%  '

# Row 18: BitcoinCash
Set-CellText $ws.Range("D18") '209.74'
Set-CellText $ws.Range("E18") '  -0.13%  '

# Row 19: Dai
Set-CellText $ws.Range("E19") '  +0.08%  '

# Row 20: Chainlink
Set-CellText $ws.Range("D20") '7.16'
Set-CellText $ws.Range("E20") '  -0.31%  '

# Row 21: Uniswap
Set-CellText $ws.Range("D21") '4.31'
Set-CellText $ws.Range("E21") '  +0.73%  '

# Row 22: Toncoin
Set-CellText $ws.Range("D22") '2.25'
Set-CellText $ws.Range("E22") '  -2.19%  '

# Row 23: Avalanche
Set-CellText $ws.Range("D23") '9.04'
Set-CellText $ws.Range("E23") '  +1.25%  '

# Row 24: Monero
Set-CellText $ws.Range("D24") '144.12'
Set-CellText $ws.Range("E24") '  +0.64%  '

# Row 25: BinanceUSD
Set-CellText $ws.Range("D25") '1.00'
Set-CellText $ws.Range("E25") '  +0.13%  '

# Row 26: Cosmos
Set-CellText $ws.Range("D26") '7.15'
Set-CellText $ws.Range("E26") '  +0.03%  '

# Row 27: Stellar
Set-CellText $ws.Range("E27") '  -0.41%  '

# Row 28: EthereumClassic
Set-CellText $ws.Range("D28") '15.38'
Set-CellText $ws.Range("E28") '  +0.31%  '

# Row 29: Hedera
Set-CellText $ws.Range("E29") '  -2.25%  '

# Row 30: PancakeSwap
Set-CellText $ws.Range("E30") '  +0.01%  '

# Row 31: Filecoin
Set-CellText $ws.Range("D31") '3.27'
Set-CellText $ws.Range("E31") '  +0.69%  '

# Row 32: InternetComputer(DFINITY)
Set-CellText $ws.Range("E32") '  +1.10%  '

# Row 33: WEMIXToken
Set-CellText $ws.Range("D33") '1.27'
Set-CellText $ws.Range("E33") '  +17.61%  '

# Row 34: Maker
Set-CellText $ws.Range("D34") '1.282.02'
Set-CellText $ws.Range("E34") '  -0.54%  '

# Row 35: HuobiToken
Set-CellText $ws.Range("E35") '  +0.63%  '

# Row 36: LidoDAOToken
Set-CellText $ws.Range("D36") '1.49'
Set-CellText $ws.Range("E36") '  +0.04%  '

# Row 37: ImmutableX
Set-CellText $ws.Range("D37") '0.595'
Set-CellText $ws.Range("E37") '  -3.90%  '

# Row 38: VeChain
Set-CellText $ws.Range("E38") '  -1.42%  '

# Row 39: ARBITRUM
Set-CellText $ws.Range("D39") '0.828'
Set-CellText $ws.Range("E39") '  -0.13%  '

# Row 40: MXToken (was FraxShare)
Set-CellText $ws.Range("B40") 'MXToken'
Set-CellText $ws.Range("C40") 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-CellText $ws.Range("D40") '2.25'
Set-CellText $ws.Range("E40") '  +2.68%  '

# Row 41: FraxShare (was MXToken)
Set-CellText $ws.Range("B41") 'FraxShare'
Set-CellText $ws.Range("C41") 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-CellText $ws.Range("D41") '5.48'
Set-CellText $ws.Range("E41") '  +0.49%  '

# Row 42: TrustWalletToken
Set-CellText $ws.Range("D42") '0.780'
Set-CellText $ws.Range("E42") '  -0.68%  '

# Row 43: Aave
Set-CellText $ws.Range("D43") '62.81'
Set-CellText $ws.Range("E43") '  -0.75%  '

# Row 44: RocketPoolETH
Set-CellText $ws.Range("D44") '1.740.83'
Set-CellText $ws.Range("E44") '  +0.26%  '

# Row 45: Quant
Set-CellText $ws.Range("D45") '90.60'
Set-CellText $ws.Range("E45") '  -0.89%  '

# Row 46: RenderToken
Set-CellText $ws.Range("D46") '1.57'
Set-CellText $ws.Range("E46") '  -0.28%  '

# Row 47: Algorand
Set-CellText $ws.Range("E47") '  +2.35%  '

# Row 48: Cronos
Set-CellText $ws.Range("D48") '0.0512'
Set-CellText $ws.Range("E48") '  +0.74%  '

# Row 49: EnergySwap
Set-CellText $ws.Range("D49") '7.56'
Set-CellText $ws.Range("E49") '  +2.98%  '

# Row 50: USDD
Set-CellText $ws.Range("E50") '  +0.13%  '

# Row 51: Mantle
Set-CellText $ws.Range("E51") '  +1.68%  '
